# Update column C ("Förändrad") date for rows 2-103 from 2023-09-02 (45171)
# to 2023-09-03 (45172) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 103; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
